# Updates to sheet "Sheet1" for rows 2-25, columns B,C,E,F,G,H,M,N,O
# (case with 380 kV done)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("B","C","E","F","G","H","M","N","O")

$data = @"
2,0.5394813449752291,0.1625366927017353,0.1508287129605463,0.4443680307746121,0.191879912830089,0.3584317889760058,0.2734228452430898,0.8824255634761045,1.007367170672637
3,0.4707715867089632,0.1505194579993798,0.1436357464271794,0.3878228170618172,0.1896484927570796,0.3609013082473638,0.2419772903493254,0.8853392644988745,1.007562709466399
4,0.4284424417755019,0.143098095030652,0.1393530083536447,0.3531389305169483,0.1884964812801329,0.3626359357430644,0.2226973546326363,0.8875120980344775,1.008593680568467
5,0.4111585908441668,0.1400633099700457,0.1376411603286982,0.3390132514313251,0.1880816411721682,0.3633976807975472,0.2148477346714017,0.888494233228073,1.009242313954601
6,0.4082865763010659,0.1395587575402999,0.1373589198077738,0.336668177824194,0.1880160491902743,0.3635274811009026,0.213544744558547,0.8886631609937936,1.009363803966011
7,0.4282094834723296,0.1430572091505411,0.1393297868385801,0.3529483938344953,0.1884906657723349,0.3626459867774727,0.2225914629503194,0.8875249517322885,1.008601503792789
8,0.5158201356663596,0.1584021511179969,0.14832066937538,0.4248636149813478,0.1910651173718776,0.3592379451142236,0.2625747237785063,0.8833506509375724,1.007245198405769
9,0.6864649748427496,0.1881459152665457,0.167024858409377,0.5661985755041457,0.1978548457477629,0.3542887756286603,0.3412003432844628,0.8782032558478079,1.011839681336738
10,0.8110890242303412,0.2097772520462797,0.1814386668591155,0.6702781546542269,0.2039205358298943,0.3517122476756072,0.3991034000233498,0.8762655277696751,1.019676722316007
11,0.867613000214817,0.2195679480328181,0.1881456479608232,0.7176906081379002,0.2069173181608477,0.3507707773495952,0.4254758266996674,0.8757828350628358,1.024219291377108
12,0.8889920290054079,0.2232681080954251,0.1907072668892695,0.7356546913071611,0.2080865256386062,0.3504474707413721,0.4354669405881282,0.8756572591509126,1.02608064054823
13,0.884388820195511,0.2224715438279077,0.1901546014916136,0.7317853510981394,0.2078331824790496,0.3505156229318089,0.4333149801351226,0.8756817618899504,1.025673476992978
14,0.8693723811657605,0.2198725114143372,0.1883559549321205,0.7191683204515869,0.2070128189275664,0.3507435129175605,0.4262977134106478,0.8757713580605326,1.024369592762895
15,0.8601710387105754,0.2182795627475969,0.1872570824428337,0.7114413442032514,0.206514808489743,0.3508874281413767,0.4220000095125869,0.8758336846824761,1.023589329744567
16,0.807391563856811,0.2091363898651366,0.1810033910881117,0.6671810134426437,0.2037294883729572,0.3517784193627449,0.3973805310191807,0.8763050858022154,1.01939957674773
17,0.7749691507852958,0.2035144948329446,0.1772055686330276,0.6400460337125793,0.202081779824951,0.3523841106192691,0.382285370833209,0.8766963164382133,1.017080091545324
18,0.7563048826119143,0.200276284195013,0.1750352716426775,0.6244449056556647,0.2011563987462921,0.3527541910513889,0.3736060587882477,0.8769588857580572,1.015837951315859
19,0.7499828129213597,0.1991790913531304,0.1743028621278739,0.619163680173358,0.2008469096327161,0.3528832199770591,0.3706679172040737,0.8770542396245133,1.015433160360828
20,0.7784222111389454,0.2041134377052458,0.1776083911615771,0.642933953830422,0.202254867620141,0.3523173872476661,0.3838919644933227,0.8766507846248857,1.017317482039402
21,0.8737837714902525,0.2206361120406939,0.1888836665443279,0.7228739723491628,0.2072528442937056,0.3506756744812094,0.4283587352876879,0.8757434899110024,1.024748738748428
22,0.9359596099488954,0.2313915986981385,0.1963800736590073,0.7751780083420101,0.2107198802167431,0.3497963023619519,0.4574461947726007,0.8754839201223348,1.030428666535443
23,0.9027891732818034,0.2256552115428008,0.1923673633892733,0.7472568307830727,0.2088510257883343,0.3502479115202277,0.4419193569464568,0.8755919951404536,1.027321658070861
24,0.7768611578147784,0.2038426746418054,0.1774262342442725,0.6416283278902171,0.2021765464029954,0.3523474847927304,0.3831656263048089,0.8766712522823212,1.017209873255069
25,0.6404294893822566,0.1801376383576496,0.1618481726750787,0.5279251897347166,0.1958299982349203,0.3554417590519989,0.3199062617507593,0.8792713493959923,1.025673476992978
"@

$lines = $data -split "`n"

foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $value = [double]$parts[$i + 1]
        $ws.Range("$colLetter$row").Value = $value
    }
}
